$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 13 (hotspring_meditation): version + localized text changes
$ws.Cells.Item(13, 2).Value = "EA 23.242"
$ws.Cells.Item(13, 3).Value = "在温泉中休息，可以温暖身心。"
$ws.Cells.Item(13, 4).Value = 'Take a "Rest" at the hot spring to warm you from the inside out!'
$ws.Cells.Item(13, 5).Value = "温泉で「休憩」すれば、心も体もほっこり"

# Add new row 14 (pond_carp)
$ws.Cells.Item(14, 1).Value = "pond_carp"
$ws.Cells.Item(14, 2).Value = "EA 23.232"
$ws.Cells.Item(14, 3).Value = "鲤鱼池"
$ws.Cells.Item(14, 4).Value = "The Carp Pond "
$ws.Cells.Item(14, 5).Value = "鯉の池"
